# Weekly update: insert a new row of "Haba" price data at the top of the
# historical series (row 65), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65 (shifts old rows 65..122 -> 66..123).
$ws.Rows("65:65").Insert()

# Populate the new row 65 with the latest observation.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44907
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112026
$ws.Range("G65").Value = "Haba"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 80
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 18000
$ws.Range("M65").Value = 18000
$ws.Range("N65").Value = "$/saco 25 kilos"
$ws.Range("O65").Value = "Región de La Araucanía"
$ws.Range("P65").Value = 720
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"
